$d = $word.ActiveDocument

# docDefaults (w:rPrDefault/w:rPr/w:rFonts) eastAsia: DejaVu Sans -> Tahoma.
# Word's object model has no direct handle onto <w:docDefaults>; the closest
# reachable surface is the per-style Font, so we apply the same eastAsia
# font change everywhere it is exposed (Normal + Heading below).

# Normal style: eastAsia DejaVu Sans -> Tahoma
$d.Styles("Normal").Font.NameFarEast = "Tahoma"

# Heading style: eastAsia DejaVu Sans -> Tahoma
$d.Styles("Heading").Font.NameFarEast = "Tahoma"

# List, Caption, Index styles: add explicit complex-script font (w:cs="DejaVu Sans")
$d.Styles("List").Font.NameBi = "DejaVu Sans"
$d.Styles("Caption").Font.NameBi = "DejaVu Sans"
$d.Styles("Index").Font.NameBi = "DejaVu Sans"

Write-Output "styles updated"
